$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the new data row (row 35) for the "Microsoft Internet Explorer Premium Edition" entry ---
$ws.Range("A35").Value2 = 34
$ws.Range("B35").Value = "Microsoft Internet Explorer Premium Edition"
$ws.Range("C35").Value = "Software"
$ws.Range("D35").Value2 = 1994
$ws.Range("E35").Value = "Y2K-2000"
$ws.Range("F35").Value2 = 36526

# Match the date formatting already used on column F (e.g. cell F34) by copying its format.
$ws.Range("F34").Copy()
$ws.Range("F35").PasteSpecial(-4122) | Out-Null
$ws.Application.CutCopyMode = $false

$ws.Range("G35").Value2 = 100
$ws.Range("H35").Value2 = 0.01
$ws.Range("I35").Value2 = 1
$ws.Range("J35").Value2 = 0
$ws.Range("K35").Value = "Internet Explorer: Still downloading; finally gave up."
$ws.Range("L35").Value = $true

# --- Update L34 (IsArchived flipped to False) and center-align it (introduces a new center-aligned cell style) ---
$ws.Range("L34").Value = $false
$ws.Range("L34").HorizontalAlignment = -4108

# --- Column width adjustments to fit the new, longer content ---
$ws.Columns.Item(2).ColumnWidth = 37.25
$ws.Columns.Item(6).ColumnWidth = 22.25
$ws.Columns.Item(7).ColumnWidth = 16.92

# --- Update the view's selection and scroll position ---
$win = $excel.ActiveWindow
$win.ScrollColumn = 3
$ws.Range("L34").Select()

# --- Page setup (orientation) ---
$ws.PageSetup.Orientation = 1

Write-Host "done"
